# Automatische test-sync: 2025-08-19 19:34:50
#
# Appends a new log entry (row 5) to the "Logs" sheet, extends the
# conditional-formatting ranges that cover the data rows so they include
# the new row, and bumps the rollup count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New row of data (mirrors the columns used by the existing rows; C and E
# are intentionally left blank, as in the source rows).
$ws.Range("A5").Value = "Demo inplannen"
$ws.Range("B5").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("D5").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("F5").Value = "2025-08-19 19:33:57"
$ws.Range("G5").Value = "Nee"
$ws.Range("H5").Value = "Ja"
$ws.Range("I5").Value = "Nee"
$ws.Range("J5").Value = "Nee"

# Grow the existing conditional-formatting rules (D/G/H/I/J, rows 2:4) so
# they cover the newly-added row 5 as well.
$ws.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D5"))
$ws.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G5"))
$ws.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H5"))
$ws.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I5"))
$ws.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J5"))

# Update the Dashboard rollup count for this category (3 -> 4).
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B2").Value = 4
